# Small update on EDA slide
# Update the "Feature Importance" row of the summary table on the
# Exploratory Data Analysis slide (slide 4) so the description reads
# "Many features have 0 importance." instead of the old hard-coded count.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$tblShape = $s.Shapes.Item(2)
$tbl = $tblShape.Table
$cell = $tbl.Cell(5, 2)
$cell.Shape.TextFrame.TextRange.Text = "Many features have 0 importance."
